# "Fixed the run mode issues."
#
# The ZohoCreateAccountTest row in the test_suite sheet had its Runmode
# flag set to "N" (skip), which was a bug - it should be "Y" (run) like
# all the other test cases. Also restore the workbook's active sheet /
# selection back to the LoginTest sheet (cell D3) instead of leaving the
# test_suite sheet (cell B5) selected.

$wb = $excel.ActiveWorkbook

# Flip the Runmode for ZohoCreateAccountTest from "N" to "Y".
$testSuite = $wb.Worksheets.Item("test_suite")
$testSuite.Range("B5").Value = "Y"

# Make LoginTest the active sheet again, with D3 selected.
$loginTest = $wb.Worksheets.Item("LoginTest")
$loginTest.Activate() | Out-Null
$loginTest.Range("D3").Select() | Out-Null
